# Spring2015.xlsx edit: "Made creation of DB/Tables if not exists, fleshed out
# new term/upload data page of test.html"
#
# Concrete spreadsheet changes:
#  1. Variables sheet ("Variables"): remove the "Year" (2015) and "Term"
#     (Spring) rows -- that data moved elsewhere -- shifting the remaining
#     Layout_* rows up.
#  2. Companies sheet ("Companies"): a new company row is inserted just
#     before the "Ryobi Die Casting, USA" row:
#       Rich-Husbands-In-Training (RHIT) / Electrical Engineering, Computer
#       Science / US Citizen,H-1 Visa / Internship/Externship, Co-op / 56
#  3. The active sheet moves from "Categories" back to "Variables".

$wb = $excel.ActiveWorkbook

$wsVariables = $wb.Worksheets.Item("Variables")
$wsCompanies = $wb.Worksheets.Item("Companies")

# --- 1. Variables: drop the Year / Term rows (rows 2 and 3) -----------------
# This removes the "Year"/"Term"/"Spring" shared strings (they become
# unreferenced) and shifts the Layout_Section* rows up to rows 2-6.
$wsVariables.Rows("2:3").Delete()

# --- 2. Companies: insert the new RHIT row before "Ryobi Die Casting, USA" --
# "Ryobi Die Casting, USA" currently lives on row 49; insert a blank row
# there and populate it, pushing Ryobi (and everything after) down to 50.
$wsCompanies.Rows("49:49").Insert()
$wsCompanies.Range("A49").Value = "Rich-Husbands-In-Training (RHIT)"
$wsCompanies.Range("B49").Value = "Electrical Engineering, Computer Science"
$wsCompanies.Range("C49").Value = "US Citizen,H-1 Visa"
$wsCompanies.Range("D49").Value = "Internship/Externship, Co-op"
$wsCompanies.Range("E49").Value = 56

# --- 3. Make "Variables" the active/selected sheet again --------------------
$wsVariables.Activate()
